$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing header cell (H1) onto the two
# new header cells so they pick up the bold/border/alignment style (s="1").
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for row 2 (plain numeric cells, no special style)
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8
